$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (time): -3 .. 21 for rows 2..26
$colA = @(-3,-2,-1,0,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21)
# Column B (tp): cycles 0..10 twice, then 1,2,3
$colB = @(0,1,2,3,4,5,6,7,8,9,10,0,1,2,3,4,5,6,7,8,9,10,1,2,3)

for ($i = 0; $i -lt $colA.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $colA[$i]
    $ws.Cells.Item($row, 2).Value = $colB[$i]
}

$ws.Range("B27").Select()
